$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29
$ws.Range("D29").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E29").Value = "['SoftwareFault']"

# Row 39
$ws.Range("D39").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E39").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"

# Row 41
$ws.Range("D41").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E41").Value = "['Normal']"
